# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new blank column immediately
# before column N ("Late"), which pushes the existing "Late" / "heading" /
# "Outstanding" columns from N/O/P to O/P/Q (and inherits the width/style of
# the column to its left, column M - "In Advance").
#
# Also make the "Repayment schedule" sheet the active sheet/tab (it was
# previously "Transactions").

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column before column N (14th column). This shifts
# "Late"/"heading"/"Outstanding" from N/O/P to O/P/Q.
$wsRepay.Columns.Item(14).Insert()

# The new column inherits the column-width of its left neighbour (M).
$wsRepay.Columns.Item(14).ColumnWidth = $wsRepay.Columns.Item(13).ColumnWidth

# The "Repayment schedule" sheet becomes the active sheet/tab (previously
# "Transactions" was active/selected).
$wsRepay.Activate()
$wsRepay.Range("R9").Select()
